$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 86.6
$ws.Cells.Item(2, 3).Value = 94.65
$ws.Cells.Item(2, 4).Value = 113.27
$ws.Cells.Item(2, 5).Value = 55.44
$ws.Cells.Item(2, 6).Value = 69.2
$ws.Cells.Item(2, 7).Value = 68.05
$ws.Cells.Item(2, 8).Value = 81.07
$ws.Cells.Item(2, 9).Value = 39.92
$ws.Cells.Item(2, 10).Value = 93.31

$ws.Cells.Item(3, 2).Value = 76.49
$ws.Cells.Item(3, 3).Value = 93.66
$ws.Cells.Item(3, 4).Value = 133.72
$ws.Cells.Item(3, 5).Value = 34.77
$ws.Cells.Item(3, 6).Value = 47.72
$ws.Cells.Item(3, 7).Value = 67.19
$ws.Cells.Item(3, 8).Value = 75.79
$ws.Cells.Item(3, 9).Value = 44.29
$ws.Cells.Item(3, 10).Value = 97.46

$ws.Cells.Item(4, 2).Value = 92.35
$ws.Cells.Item(4, 3).Value = 80.51
$ws.Cells.Item(4, 4).Value = 129.72
$ws.Cells.Item(4, 5).Value = 40.65
$ws.Cells.Item(4, 6).Value = 33.33
$ws.Cells.Item(4, 7).Value = 59.77
$ws.Cells.Item(4, 8).Value = 82.69
$ws.Cells.Item(4, 9).Value = 36.93
$ws.Cells.Item(4, 10).Value = 69.98

$ws.Cells.Item(5, 2).Value = 70.01
$ws.Cells.Item(5, 3).Value = 100.56
$ws.Cells.Item(5, 4).Value = 76.56
$ws.Cells.Item(5, 5).Value = 43.09
$ws.Cells.Item(5, 6).Value = 82.05
$ws.Cells.Item(5, 7).Value = 29.63
$ws.Cells.Item(5, 8).Value = 65.52
$ws.Cells.Item(5, 9).Value = 75.69
$ws.Cells.Item(5, 10).Value = 83.86

$ws.Cells.Item(6, 2).Value = 107.27
$ws.Cells.Item(6, 3).Value = 161.29
$ws.Cells.Item(6, 4).Value = 201.45
$ws.Cells.Item(6, 5).Value = 41.46
$ws.Cells.Item(6, 6).Value = 47.06
$ws.Cells.Item(6, 7).Value = 60.32
$ws.Cells.Item(6, 8).Value = 94.13
$ws.Cells.Item(6, 9).Value = 42.62
$ws.Cells.Item(6, 10).Value = 135.83

$ws.Cells.Item(7, 2).Value = 64.46
$ws.Cells.Item(7, 3).Value = 124.07
$ws.Cells.Item(7, 4).Value = 181.91
$ws.Cells.Item(7, 5).Value = 68.29
$ws.Cells.Item(7, 6).Value = 75.0
$ws.Cells.Item(7, 7).Value = 142.86
$ws.Cells.Item(7, 8).Value = 69.18
$ws.Cells.Item(7, 9).Value = 30.87
$ws.Cells.Item(7, 10).Value = 110.53

$ws.Cells.Item(8, 2).Value = 66.88
$ws.Cells.Item(8, 3).Value = 50.43
$ws.Cells.Item(8, 4).Value = 132.98
$ws.Cells.Item(8, 5).Value = 4.42
$ws.Cells.Item(8, 6).Value = 14.04
$ws.Cells.Item(8, 7).Value = 78.33
$ws.Cells.Item(8, 8).Value = 70.42
$ws.Cells.Item(8, 9).Value = 30.13
$ws.Cells.Item(8, 10).Value = 99.81

$ws.Cells.Item(9, 2).Value = 67.34
$ws.Cells.Item(9, 3).Value = 69.71
$ws.Cells.Item(9, 4).Value = 86.02
$ws.Cells.Item(9, 5).Value = 13.82
$ws.Cells.Item(9, 6).Value = 47.62
$ws.Cells.Item(9, 7).Value = 41.67
$ws.Cells.Item(9, 8).Value = 70.98
$ws.Cells.Item(9, 9).Value = 73.77
$ws.Cells.Item(9, 10).Value = 79.81

$ws.Cells.Item(10, 2).Value = 90.43
$ws.Cells.Item(10, 3).Value = 84.65
$ws.Cells.Item(10, 4).Value = 92.81
$ws.Cells.Item(10, 5).Value = 30.08
$ws.Cells.Item(10, 6).Value = 91.67
$ws.Cells.Item(10, 7).Value = 77.14
$ws.Cells.Item(10, 8).Value = 88.02
$ws.Cells.Item(10, 9).Value = 33.99
$ws.Cells.Item(10, 10).Value = 90.55

$ws.Cells.Item(11, 2).Value = 63.26
$ws.Cells.Item(11, 3).Value = 53.47
$ws.Cells.Item(11, 4).Value = 41.88
$ws.Cells.Item(11, 5).Value = 9.98
$ws.Cells.Item(11, 6).Value = 29.63
$ws.Cells.Item(11, 7).Value = 53.33
$ws.Cells.Item(11, 8).Value = 96.67
$ws.Cells.Item(11, 9).Value = 20.9
$ws.Cells.Item(11, 10).Value = 53.69

$ws.Cells.Item(12, 2).Value = 109.72
$ws.Cells.Item(12, 3).Value = 92.24
$ws.Cells.Item(12, 4).Value = 109.18
$ws.Cells.Item(12, 5).Value = 72.93
$ws.Cells.Item(12, 6).Value = 145.45
$ws.Cells.Item(12, 7).Value = 101.27
$ws.Cells.Item(12, 8).Value = 81.42
$ws.Cells.Item(12, 9).Value = 43.88
$ws.Cells.Item(12, 10).Value = 97.78

$ws.Cells.Item(13, 2).Value = 74.43
$ws.Cells.Item(13, 3).Value = 86.23
$ws.Cells.Item(13, 4).Value = 128.39
$ws.Cells.Item(13, 5).Value = 22.76
$ws.Cells.Item(13, 6).Value = 106.67
$ws.Cells.Item(13, 7).Value = 86.67
$ws.Cells.Item(13, 8).Value = 102.18
$ws.Cells.Item(13, 9).Value = 49.82
$ws.Cells.Item(13, 10).Value = 116.67

$ws.Cells.Item(14, 2).Value = 104.76
$ws.Cells.Item(14, 3).Value = 138.87
$ws.Cells.Item(14, 4).Value = 113.63
$ws.Cells.Item(14, 5).Value = 4.88
$ws.Cells.Item(14, 6).Value = 248.89
$ws.Cells.Item(14, 7).Value = 80.95
$ws.Cells.Item(14, 8).Value = 131.5
$ws.Cells.Item(14, 9).Value = 15.32
$ws.Cells.Item(14, 10).Value = 94.12

$ws.Cells.Item(15, 2).Value = 52.12
$ws.Cells.Item(15, 3).Value = 35.22
$ws.Cells.Item(15, 4).Value = 78.89
$ws.Cells.Item(15, 5).Value = 0.0
$ws.Cells.Item(15, 6).Value = 0.0
$ws.Cells.Item(15, 7).Value = 54.17
$ws.Cells.Item(15, 8).Value = 60.45
$ws.Cells.Item(15, 9).Value = 46.58
$ws.Cells.Item(15, 10).Value = 67.58

$ws.Cells.Item(16, 2).Value = 123.86
$ws.Cells.Item(16, 3).Value = 92.07
$ws.Cells.Item(16, 4).Value = 113.93
$ws.Cells.Item(16, 5).Value = 63.41
$ws.Cells.Item(16, 6).Value = 9.52
$ws.Cells.Item(16, 7).Value = 58.48
$ws.Cells.Item(16, 8).Value = 54.98
$ws.Cells.Item(16, 9).Value = 34.5
$ws.Cells.Item(16, 10).Value = 120.0

$ws.Cells.Item(17, 2).Value = 106.8
$ws.Cells.Item(17, 3).Value = 121.12
$ws.Cells.Item(17, 4).Value = 120.15
$ws.Cells.Item(17, 5).Value = 79.68
$ws.Cells.Item(17, 6).Value = 92.1
$ws.Cells.Item(17, 7).Value = 81.25
$ws.Cells.Item(17, 8).Value = 83.33
$ws.Cells.Item(17, 9).Value = 41.95
$ws.Cells.Item(17, 10).Value = 96.87

$ws.Cells.Item(18, 2).Value = 111.85
$ws.Cells.Item(18, 3).Value = 102.93
$ws.Cells.Item(18, 4).Value = 118.64
$ws.Cells.Item(18, 5).Value = 145.86
$ws.Cells.Item(18, 6).Value = 133.33
$ws.Cells.Item(18, 7).Value = 84.85
$ws.Cells.Item(18, 8).Value = 96.86
$ws.Cells.Item(18, 9).Value = 61.02
$ws.Cells.Item(18, 10).Value = 90.93

$ws.Cells.Item(19, 2).Value = 109.95
$ws.Cells.Item(19, 3).Value = 120.75
$ws.Cells.Item(19, 4).Value = 110.43
$ws.Cells.Item(19, 5).Value = 103.87
$ws.Cells.Item(19, 6).Value = 85.71
$ws.Cells.Item(19, 7).Value = 72.38
$ws.Cells.Item(19, 8).Value = 68.24
$ws.Cells.Item(19, 9).Value = 34.12
$ws.Cells.Item(19, 10).Value = 105.88

$ws.Cells.Item(20, 2).Value = 114.57
$ws.Cells.Item(20, 3).Value = 148.27
$ws.Cells.Item(20, 4).Value = 120.17
$ws.Cells.Item(20, 5).Value = 85.52
$ws.Cells.Item(20, 6).Value = 5.8
$ws.Cells.Item(20, 7).Value = 81.48
$ws.Cells.Item(20, 8).Value = 118.09
$ws.Cells.Item(20, 9).Value = 36.63
$ws.Cells.Item(20, 10).Value = 94.34

$ws.Cells.Item(21, 2).Value = 108.06
$ws.Cells.Item(21, 3).Value = 146.05
$ws.Cells.Item(21, 4).Value = 122.8
$ws.Cells.Item(21, 5).Value = 22.1
$ws.Cells.Item(21, 6).Value = 278.79
$ws.Cells.Item(21, 7).Value = 87.88
$ws.Cells.Item(21, 8).Value = 81.65
$ws.Cells.Item(21, 9).Value = 44.87
$ws.Cells.Item(21, 10).Value = 89.46

$ws.Cells.Item(22, 2).Value = 82.47
$ws.Cells.Item(22, 3).Value = 100.0
$ws.Cells.Item(22, 4).Value = 131.65
$ws.Cells.Item(22, 5).Value = 40.52
$ws.Cells.Item(22, 6).Value = 50.0
$ws.Cells.Item(22, 7).Value = 78.79
$ws.Cells.Item(22, 8).Value = 70.05
$ws.Cells.Item(22, 9).Value = 35.6
$ws.Cells.Item(22, 10).Value = 102.7

$ws.Cells.Item(23, 2).Value = 84.12
$ws.Cells.Item(23, 3).Value = 93.79
$ws.Cells.Item(23, 4).Value = 106.54
$ws.Cells.Item(23, 5).Value = 31.14
$ws.Cells.Item(23, 6).Value = 51.85
$ws.Cells.Item(23, 7).Value = 53.91
$ws.Cells.Item(23, 8).Value = 77.39
$ws.Cells.Item(23, 9).Value = 36.52
$ws.Cells.Item(23, 10).Value = 88.39

$ws.Cells.Item(24, 2).Value = 37.46
$ws.Cells.Item(24, 3).Value = 74.64
$ws.Cells.Item(24, 4).Value = 51.35
$ws.Cells.Item(24, 5).Value = 12.52
$ws.Cells.Item(24, 6).Value = 0.0
$ws.Cells.Item(24, 7).Value = 23.7
$ws.Cells.Item(24, 8).Value = 25.41
$ws.Cells.Item(24, 9).Value = 21.23
$ws.Cells.Item(24, 10).Value = 37.5

$ws.Cells.Item(25, 2).Value = 133.06
$ws.Cells.Item(25, 3).Value = 110.51
$ws.Cells.Item(25, 4).Value = 193.19
$ws.Cells.Item(25, 5).Value = 36.59
$ws.Cells.Item(25, 6).Value = 133.33
$ws.Cells.Item(25, 7).Value = 87.18
$ws.Cells.Item(25, 8).Value = 95.22
$ws.Cells.Item(25, 9).Value = 43.62
$ws.Cells.Item(25, 10).Value = 137.39

$ws.Cells.Item(26, 2).Value = 82.25
$ws.Cells.Item(26, 3).Value = 79.22
$ws.Cells.Item(26, 4).Value = 124.55
$ws.Cells.Item(26, 5).Value = 37.4
$ws.Cells.Item(26, 6).Value = 66.67
$ws.Cells.Item(26, 7).Value = 32.52
$ws.Cells.Item(26, 8).Value = 73.44
$ws.Cells.Item(26, 9).Value = 50.0
$ws.Cells.Item(26, 10).Value = 114.29

$ws.Cells.Item(27, 2).Value = 82.92
$ws.Cells.Item(27, 3).Value = 142.02
$ws.Cells.Item(27, 4).Value = 88.62
$ws.Cells.Item(27, 5).Value = 40.22
$ws.Cells.Item(27, 6).Value = 31.37
$ws.Cells.Item(27, 7).Value = 49.23
$ws.Cells.Item(27, 8).Value = 96.62
$ws.Cells.Item(27, 9).Value = 24.92
$ws.Cells.Item(27, 10).Value = 68.33

$ws.Cells.Item(28, 2).Value = 50.8
$ws.Cells.Item(28, 3).Value = 49.01
$ws.Cells.Item(28, 4).Value = 77.19
$ws.Cells.Item(28, 5).Value = 37.4
$ws.Cells.Item(28, 6).Value = 38.89
$ws.Cells.Item(28, 7).Value = 33.33
$ws.Cells.Item(28, 8).Value = 82.19
$ws.Cells.Item(28, 9).Value = 25.87
$ws.Cells.Item(28, 10).Value = 90.0

$ws.Cells.Item(29, 2).Value = 81.79
$ws.Cells.Item(29, 3).Value = 72.18
$ws.Cells.Item(29, 4).Value = 94.51
$ws.Cells.Item(29, 5).Value = 8.84
$ws.Cells.Item(29, 6).Value = 24.24
$ws.Cells.Item(29, 7).Value = 63.49
$ws.Cells.Item(29, 8).Value = 79.48
$ws.Cells.Item(29, 9).Value = 44.95
$ws.Cells.Item(29, 10).Value = 78.38

$ws.Cells.Item(30, 2).Value = 77.69
$ws.Cells.Item(30, 3).Value = 86.7
$ws.Cells.Item(30, 4).Value = 118.0
$ws.Cells.Item(30, 5).Value = 113.35
$ws.Cells.Item(30, 6).Value = 64.31
$ws.Cells.Item(30, 7).Value = 67.68
$ws.Cells.Item(30, 8).Value = 79.27
$ws.Cells.Item(30, 9).Value = 42.57
$ws.Cells.Item(30, 10).Value = 92.64

$ws.Cells.Item(31, 2).Value = 97.12
$ws.Cells.Item(31, 3).Value = 137.94
$ws.Cells.Item(31, 4).Value = 148.34
$ws.Cells.Item(31, 5).Value = 31.71
$ws.Cells.Item(31, 6).Value = 25.0
$ws.Cells.Item(31, 7).Value = 65.61
$ws.Cells.Item(31, 8).Value = 94.24
$ws.Cells.Item(31, 9).Value = 58.55
$ws.Cells.Item(31, 10).Value = 106.67

$ws.Cells.Item(32, 2).Value = 66.09
$ws.Cells.Item(32, 3).Value = 133.33
$ws.Cells.Item(32, 4).Value = 126.61
$ws.Cells.Item(32, 5).Value = 40.65
$ws.Cells.Item(32, 6).Value = 53.33
$ws.Cells.Item(32, 7).Value = 60.22
$ws.Cells.Item(32, 8).Value = 75.38
$ws.Cells.Item(32, 9).Value = 37.81
$ws.Cells.Item(32, 10).Value = 77.45

$ws.Cells.Item(33, 2).Value = 75.26
$ws.Cells.Item(33, 3).Value = 69.82
$ws.Cells.Item(33, 4).Value = 129.72
$ws.Cells.Item(33, 5).Value = 353.66
$ws.Cells.Item(33, 6).Value = 31.37
$ws.Cells.Item(33, 7).Value = 83.72
$ws.Cells.Item(33, 8).Value = 52.7
$ws.Cells.Item(33, 9).Value = 46.71
$ws.Cells.Item(33, 10).Value = 66.67

$ws.Cells.Item(34, 2).Value = 67.8
$ws.Cells.Item(34, 3).Value = 50.29
$ws.Cells.Item(34, 4).Value = 96.49
$ws.Cells.Item(34, 5).Value = 125.2
$ws.Cells.Item(34, 6).Value = 13.33
$ws.Cells.Item(34, 7).Value = 35.9
$ws.Cells.Item(34, 8).Value = 89.62
$ws.Cells.Item(34, 9).Value = 57.93
$ws.Cells.Item(34, 10).Value = 112.5

$ws.Cells.Item(35, 2).Value = 85.75
$ws.Cells.Item(35, 3).Value = 95.58
$ws.Cells.Item(35, 4).Value = 99.8
$ws.Cells.Item(35, 5).Value = 95.12
$ws.Cells.Item(35, 6).Value = 195.56
$ws.Cells.Item(35, 7).Value = 73.33
$ws.Cells.Item(35, 8).Value = 81.72
$ws.Cells.Item(35, 9).Value = 30.26
$ws.Cells.Item(35, 10).Value = 105.15

$ws.Cells.Item(36, 2).Value = 59.17
$ws.Cells.Item(36, 3).Value = 37.48
$ws.Cells.Item(36, 4).Value = 83.06
$ws.Cells.Item(36, 5).Value = 41.25
$ws.Cells.Item(36, 6).Value = 54.9
$ws.Cells.Item(36, 7).Value = 84.44
$ws.Cells.Item(36, 8).Value = 86.07
$ws.Cells.Item(36, 9).Value = 30.26
$ws.Cells.Item(36, 10).Value = 99.35
